$d = $word.ActiveDocument

function Set-ParagraphXml($range, [string]$innerParagraphXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerParagraphXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$range.InsertXML($pkg)
}

# --- Paragraph "Gestionar Publicación: ..." -> "Gestionar Publicaciones: ..." ---
# Find the paragraph by its distinctive leading text so we grab the whole paragraph range.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Gestionar Publicaci*Permite*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newPara1 = '<w:p w:rsidR="00044E6A" w:rsidRPr="00044E6A" w:rsidRDefault="008A48C3" w:rsidP="00016624">' +
        '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r w:rsidRPr="0000359D"><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>Gestionar</w:t></w:r>' +
        '<w:r w:rsidR="00044E6A" w:rsidRPr="0000359D"><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r w:rsidRPr="0000359D"><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>Publicacio</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>n</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>es</w:t></w:r>' +
        '<w:r w:rsidR="00044E6A"><w:t>: Permite</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> al usuario crear, modificar o eliminar una publicación en su mur</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t>o o en el muro de los amigos en los que tenga permiso.</w:t></w:r>' +
        '</w:p>'
    Set-ParagraphXml $target.Range $newPara1
}

# --- Paragraph "Enviar Notificaciones: ..." -> drop the stale _GoBack bookmark ---
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Enviar Notificaciones*") {
        $target2 = $p
        break
    }
}

if ($target2 -ne $null) {
    $newPara2 = '<w:p w:rsidR="00E8253F" w:rsidRPr="00E8253F" w:rsidRDefault="00E8253F" w:rsidP="00016624">' +
        '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t>Enviar Notificaciones</w:t></w:r>' +
        '<w:r><w:t>: El sistema envía notificaciones o mail.</w:t></w:r>' +
        '</w:p>'
    Set-ParagraphXml $target2.Range $newPara2
}
